# Adds 5 new OAUTH negative-test rows (rows 6-10) to the OAUTH worksheet,
# together with a hyperlink on G8, matching the "Added testcases in OAUTH" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 6: wrong username ----
$ws.Range("A6").Value = "OPQA-XXXX5"
$ws.Range("B6").Value = "Verify that to get error status by passing wrong username in OAUTH API"
$ws.Range("C6").Value = "1POAUTH"
$ws.Range("D6").Value = "/token"
$ws.Range("E6").Value = "POST"
$ws.Range("F6").Value = "Content-Type=application/x-www-form-urlencoded||Authorization=Basic YWRtaW46cGFzc3dvcmQ="
$ws.Range("G6").Value = "?grant_type=password&username=(ddMMMyyyy_HHmmss)1@tr.com&password=Neon@123"
$ws.Range("J6").Value = "status=400||code=40020||error_description=Login failed because user profile does not match||error=invalid_grant"
$ws.Rows.Item(6).RowHeight = 90

# ---- Row 7: wrong password ----
$ws.Range("A7").Value = "OPQA-XXXX6"
$ws.Range("B7").Value = "Verify that to get error status by passing wrong password in OAUTH API"
$ws.Range("C7").Value = "1POAUTH"
$ws.Range("D7").Value = "/token"
$ws.Range("E7").Value = "POST"
$ws.Range("F7").Value = "Content-Type=application/x-www-form-urlencoded||Authorization=Basic YWRtaW46cGFzc3dvcmQ="
$ws.Range("G7").Value = "?grant_type=password&username=(ddMMMyyyy_HHmmss)@tr.com&password=Neon@1234"
$ws.Range("J7").Value = "status=400||code=40012||error_description=Login failed||error=invalid_grant"
$ws.Rows.Item(7).RowHeight = 60

# ---- Row 8: empty username (with a hyperlink on the querystring cell) ----
$ws.Range("A8").Value = "OPQA-XXXX7"
$ws.Range("B8").Value = "Verify that to get error status by passing empty username in OAUTH API"
$ws.Range("C8").Value = "1POAUTH"
$ws.Range("D8").Value = "/token"
$ws.Range("E8").Value = "POST"
$ws.Range("F8").Value = "Content-Type=application/x-www-form-urlencoded||Authorization=Basic YWRtaW46cGFzc3dvcmQ="
$ws.Range("G8").Value = "?grant_type=password&username=&password=Neon@123"
$ws.Range("J8").Value = "status=400||error_description=Missing parameters: username||error=invalid_request"
$ws.Rows.Item(8).RowHeight = 60

# ---- Row 9: empty password ----
$ws.Range("A9").Value = "OPQA-XXXX8"
$ws.Range("B9").Value = "Verify that to get error status by passing empty password in OAUTH API"
$ws.Range("C9").Value = "1POAUTH"
$ws.Range("D9").Value = "/token"
$ws.Range("E9").Value = "POST"
$ws.Range("F9").Value = "Content-Type=application/x-www-form-urlencoded||Authorization=Basic YWRtaW46cGFzc3dvcmQ="
$ws.Range("G9").Value = "?grant_type=password&username=(ddMMMyyyy_HHmmss)1@tr.com&password="
$ws.Range("J9").Value = "status=400||error_description=Missing parameters: password||error=invalid_request"
$ws.Rows.Item(9).RowHeight = 60

# ---- Row 10: empty username and password ----
$ws.Range("A10").Value = "OPQA-XXXX9"
$ws.Range("B10").Value = "Verify that to get error status by passing empty username and password in OAUTH API"
$ws.Range("C10").Value = "1POAUTH"
$ws.Range("D10").Value = "/token"
$ws.Range("E10").Value = "POST"
$ws.Range("F10").Value = "Content-Type=application/x-www-form-urlencoded||Authorization=Basic YWRtaW46cGFzc3dvcmQ="
$ws.Range("G10").Value = "?grant_type=password&username=&password="
$ws.Range("J10").Value = "status=400||error_description=Missing parameters: password username||error=invalid_request"
$ws.Rows.Item(10).RowHeight = 75

# Hyperlink on the querystring of the "empty username" test case (G8)
$ws.Hyperlinks.Add($ws.Range("G8"), "https://jira.example.com/browse/OPQA-XXXX7")

# Keep the selection in sync with the now-larger data range (L2:L10)
[void]$ws.Range("L2:L10").Select()
